$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting from H1 (bold, centered, bordered header style) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate columns I and J (rows 2-72) with data values
$arr = New-Object 'object[,]' 71,2
$arr[0,0] = 7
$arr[0,1] = 8
$arr[1,0] = 7
$arr[1,1] = 8
$arr[2,0] = 6
$arr[2,1] = 6
$arr[3,0] = 8
$arr[3,1] = 8
$arr[4,0] = 7
$arr[4,1] = 7
$arr[5,0] = 7
$arr[5,1] = 7
$arr[6,0] = 8
$arr[6,1] = 8
$arr[7,0] = 6
$arr[7,1] = 8
$arr[8,0] = 7
$arr[8,1] = 7
$arr[9,0] = 9
$arr[9,1] = 9
$arr[10,0] = 8
$arr[10,1] = 8
$arr[11,0] = 8
$arr[11,1] = 8
$arr[12,0] = 7
$arr[12,1] = 7
$arr[13,0] = 8
$arr[13,1] = 8
$arr[14,0] = 7
$arr[14,1] = 7
$arr[15,0] = 7
$arr[15,1] = 7
$arr[16,0] = 10
$arr[16,1] = 10
$arr[17,0] = 7
$arr[17,1] = 7
$arr[18,0] = 8
$arr[18,1] = 8
$arr[19,0] = 7
$arr[19,1] = 7
$arr[20,0] = 7
$arr[20,1] = 7
$arr[21,0] = 7
$arr[21,1] = 7
$arr[22,0] = 6
$arr[22,1] = 6
$arr[23,0] = 7
$arr[23,1] = 8
$arr[24,0] = 7
$arr[24,1] = 7
$arr[25,0] = 6
$arr[25,1] = 7
$arr[26,0] = 7
$arr[26,1] = 7
$arr[27,0] = 10
$arr[27,1] = 10
$arr[28,0] = 9
$arr[28,1] = 9
$arr[29,0] = 8
$arr[29,1] = 8
$arr[30,0] = 7
$arr[30,1] = 7
$arr[31,0] = 5
$arr[31,1] = 5
$arr[32,0] = 7
$arr[32,1] = 8
$arr[33,0] = 7
$arr[33,1] = 7
$arr[34,0] = 5
$arr[34,1] = 5
$arr[35,0] = 6
$arr[35,1] = 7
$arr[36,0] = 10
$arr[36,1] = 10
$arr[37,0] = 7
$arr[37,1] = 7
$arr[38,0] = 7
$arr[38,1] = 7
$arr[39,0] = 6
$arr[39,1] = 6
$arr[40,0] = 7
$arr[40,1] = 7
$arr[41,0] = 7
$arr[41,1] = 7
$arr[42,0] = 6
$arr[42,1] = 6
$arr[43,0] = 7
$arr[43,1] = 7
$arr[44,0] = 7
$arr[44,1] = 7
$arr[45,0] = 10
$arr[45,1] = 10
$arr[46,0] = 7
$arr[46,1] = 7
$arr[47,0] = 5
$arr[47,1] = 5
$arr[48,0] = 7
$arr[48,1] = 7
$arr[49,0] = 5
$arr[49,1] = 6
$arr[50,0] = 6
$arr[50,1] = 6
$arr[51,0] = 7
$arr[51,1] = 7
$arr[52,0] = 9
$arr[52,1] = 10
$arr[53,0] = 8
$arr[53,1] = 8
$arr[54,0] = 8
$arr[54,1] = 8
$arr[55,0] = 6
$arr[55,1] = 6
$arr[56,0] = 8
$arr[56,1] = 8
$arr[57,0] = 8
$arr[57,1] = 9
$arr[58,0] = 4
$arr[58,1] = 4
$arr[59,0] = 5
$arr[59,1] = 5
$arr[60,0] = 7
$arr[60,1] = 7
$arr[61,0] = 8
$arr[61,1] = 8
$arr[62,0] = 8
$arr[62,1] = 8
$arr[63,0] = 5
$arr[63,1] = 5
$arr[64,0] = 9
$arr[64,1] = 9
$arr[65,0] = 6
$arr[65,1] = 6
$arr[66,0] = 8
$arr[66,1] = 8
$arr[67,0] = 7
$arr[67,1] = 7
$arr[68,0] = 6
$arr[68,1] = 6
$arr[69,0] = 5
$arr[69,1] = 5
$arr[70,0] = 5
$arr[70,1] = 5

$ws.Range("I2:J72").Value = $arr

Write-Host "Edit complete"